$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the refreshed view settings (100% zoom) recorded for the sheet.
$excel.ActiveWindow.Zoom = 100

# The header row previously carried a bold font + teal fill style; remove it
# so the header cells fall back to the workbook's default (unstyled) look.
$ws.Range("A1:D1").ClearFormats()

# Header text for column D changes from "string" to "STRING"
$ws.Range("D1").Value = "STRING"

# New "missing" column (B) gets populated for rows 2-4:
#  - B2 stays truly blank (touch it without giving it a value/style)
$ws.Range("B2").Borders.LineStyle = -4142   # xlLineStyleNone

#  - B3 becomes an empty string
$ws.Range("B3").Value = "'"
$ws.Range("B3").ClearFormats()

#  - B4 becomes the literal text "#NAME?" (force text so it isn't read as an error)
$ws.Range("B4").Value = "'#NAME?"
$ws.Range("B4").ClearFormats()
